# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment schedule"
#   sheet, shifting the old N/O/P columns (Late / heading / Outstanding) one
#   position to the right, and give the new column the same width as column M.
# - Make "Repayment schedule" the active sheet/tab, with I16 selected there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column before the existing column N (shifts N->O, O->P, P->Q)
$ws.Columns("N").Insert()

# New column N should match the width of column M ("In Advance")
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet, with I16 selected
$ws.Activate() | Out-Null
$ws.Range("I16").Select() | Out-Null
